# test cases and error proofing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix / update header row (row 2) ---
$ws.Range("A2").Value = "Name"
$ws.Range("D2").Value = "money you have with you"
$ws.Range("E2").Value = "flyers bought"
$ws.Range("J2").Value = "Location"
$ws.Range("K2").Value = "How game ended"

# --- Test case rows (3-10) ---

# Row 3 - Max
$ws.Range("A3").Value = "Max "
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = "no"
$ws.Range("D3").Value = "n/a"
$ws.Range("E3").Value = "n/a"
$ws.Range("F3").Value = "n/a"
$ws.Range("J3").Value = "Rockville"
$ws.Range("K3").Value = "You went home empty handed"

# Row 4 - Tim
$ws.Range("A4").Value = "Tim"
$ws.Range("B4").Value = 101
$ws.Range("C4").Value = "n/a"
$ws.Range("D4").Value = "n/a"
$ws.Range("E4").Value = "n/a"
$ws.Range("F4").Value = "n/a"
$ws.Range("J4").Value = "None"
$ws.Range("K4").Value = "That is a little far for a walk, why don’t you go somwhere closer"

# Row 5 - John
$ws.Range("A5").Value = "John"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "yes"
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = "n/a"
$ws.Range("F5").Value = "n/a"
$ws.Range("J5").Value = "Rockville"
$ws.Range("K5").Value = "you went home broke"

# Row 6 - Eric
$ws.Range("A6").Value = "Eric"
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = "n/a"
$ws.Range("D6").Value = "n/a"
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = "n/a"
$ws.Range("J6").Value = "Forest"
$ws.Range("K6").Value = "you went home unlike Jeffs cat… you should of helped"

# Row 7 - Phil
$ws.Range("A7").Value = "Phil"
$ws.Range("B7").Value = 90
$ws.Range("C7").Value = "n/a"
$ws.Range("D7").Value = "n/a"
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = "n/a"
$ws.Range("J7").Value = "Forest"
$ws.Range("K7").Value = "you went home unlike Jeffs cat… you could of printed more flyers you know"

# Row 8 - Fred
$ws.Range("A8").Value = "Fred"
$ws.Range("B8").Value = 50
$ws.Range("C8").Value = "n/a"
$ws.Range("D8").Value = "nla"
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = "gold"
$ws.Range("J8").Value = "forest"
$ws.Range("K8").Value = "you went home rich and Jeff found his cat"

# Row 9 - Eric
$ws.Range("A9").Value = "Eric"
$ws.Range("B9").Value = 60
$ws.Range("C9").Value = "n/a"
$ws.Range("D9").Value = "n/a"
$ws.Range("E9").Value = 109
$ws.Range("F9").Value = "a free ride home"
$ws.Range("J9").Value = "forest"
$ws.Range("K9").Value = "you got home safe and Jeff found his cat with your help"

# Row 10 - Lewis
$ws.Range("A10").Value = "Lewis"
$ws.Range("B10").Value = 12
$ws.Range("C10").Value = "n/a"
$ws.Range("D10").Value = "n/a"
$ws.Range("E10").Value = 67
$ws.Range("F10").Value = "Rubies"
$ws.Range("J10").Value = "forest"
$ws.Range("K10").Value = "you went home rish and Jeff found his cat"

# --- Selection matches the saved state in the workbook ---
$ws.Range("H14").Select()
